$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metrics")

$ws.Range("B2").Value = 220753.22
$ws.Range("B3").Value = 181213.44
$ws.Range("B4").Value = 70198.47
$ws.Range("B5").Value = 8761
$ws.Range("B6").Value = 4587884.6899999995
$ws.Range("B7").Value = 3871032.11
$ws.Range("B8").Value = 1340800.6100000001
$ws.Range("B9").Value = 177762
$ws.Range("B10").Value = 33053208.490999825
$ws.Range("B11").Value = 31146253.630000003
$ws.Range("B12").Value = 11622509.500000002
$ws.Range("B13").Value = 1275389

$ws.Activate()
$ws.Range("D13").Select()
